$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "260.85"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "2.03%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "27.23"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "2.95%"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.691"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "0.56%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.06122"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "3.33%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.655"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "0.79%"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8536"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "0.16%"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9207"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "1.09%"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1399"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "1.68%"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.04667"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "11.22%"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07085"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "1.35%"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03055"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "1.09%"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.28%"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001537"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1.22%"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0006083"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "0.67%"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006051"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "-0.68%"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.450"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-0.54%"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.145"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.10%"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-0.63%"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.1310"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "2.05%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.084"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "5.97%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04242"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "1.06%"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "0.11%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.003800"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "-18.93%"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "0.03%"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0001575"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "3.47%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03871"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "2.15%"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1114"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.79%"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-34.85%"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "12.45%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.002217"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-4.08%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005163"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "0.23%"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "0.06%"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "8.13%"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1622"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "-32.65%"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "0.06%"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "0.06%"
